# Update the header labels for the severity-level columns (E1:L1) on every
# worksheet in the workbook, per the commit:
#   "% 1-2" -> "% severity levels 1-2"
#   "# 1-2" -> "# severity levels 1-2"
#   "% 3"   -> "% severity level 3"
#   "# 3"   -> "# severity level 3"
#   "% 4"   -> "% severity level 4"
#   "# 4"   -> "# severity level 4"
#   "% 5"   -> "% severity level 5"
#   "# 5"   -> "# severity level 5"

$wb = $excel.ActiveWorkbook

$headerMap = @{
    "% 1-2" = "% severity levels 1-2"
    "# 1-2" = "# severity levels 1-2"
    "% 3"   = "% severity level 3"
    "# 3"   = "# severity level 3"
    "% 4"   = "% severity level 4"
    "# 4"   = "# severity level 4"
    "% 5"   = "% severity level 5"
    "# 5"   = "# severity level 5"
}

foreach ($ws in $wb.Worksheets) {
    foreach ($col in @("E", "F", "G", "H", "I", "J", "K", "L")) {
        $cell = $ws.Range($col + "1")
        $current = $cell.Value()
        if ($headerMap.ContainsKey($current)) {
            $cell.Value = $headerMap[$current]
        }
    }
}
